# Insert two new data rows right after existing row 30 (pushing former rows 31..85 down to 33..87)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31:A32").EntireRow.Insert()

# New row 31: Damasco, Dina, Especial, Terminal La Palmera de La Serena / Coquimbo
$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44914
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103003
$ws.Range("J31").Value = "Damasco"
$ws.Range("K31").Value = "Dina"
$ws.Range("L31").Value = "Especial"
$ws.Range("M31").Value = 360
$ws.Range("N31").Value = 22000
$ws.Range("O31").Value = 23000
$ws.Range("P31").Value = 22500
$ws.Range("Q31").Value = "$/caja 16 kilos"
$ws.Range("R31").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S31").Value = 1406
$ws.Range("T31").Value = 16

# New row 32: Damasco, Dina, Primera, Terminal La Palmera de La Serena / Coquimbo
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44914
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = "Frutos de hueso (carozo)"
$ws.Range("I32").Value = 100103003
$ws.Range("J32").Value = "Damasco"
$ws.Range("K32").Value = "Dina"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 19000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 19500
$ws.Range("Q32").Value = "$/caja 16 kilos"
$ws.Range("R32").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S32").Value = 1219
$ws.Range("T32").Value = 16
